# Day 0 exercises + Day 1: ex. 1 to 3
# Add description for the "nighttime_lights" column (row 34) in the data key sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the previous "description" cell (B33) so the new
# cell picks up the same style used by all the other column-description cells.
$ws.Range("B33").Copy()
$ws.Range("B34").PasteSpecial(-4122)

$ws.Range("B34").Value = "Brightness of night time lights in that location (arbitrary scale). Defined in https://pophealthmetrics.biomedcentral.com/articles/10.1186/1478-7954-6-5#Sec10"

# Restore the view/selection state to where the author left off editing.
$excel.ActiveWindow.ScrollRow = 21
$ws.Range("E26").Select() | Out-Null
